$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated meanEMG legmaxROM values (row 1 header + CON/STR rows 2-3), columns B:E
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

$ws.Range("B2").Value = 237.11241462348252
$ws.Range("C2").Value = 173.8357176670649
$ws.Range("D2").Value = 238.50370248269837
$ws.Range("E2").Value = 175.68100970774961

$ws.Range("B3").Value = 217.36520251781573
$ws.Range("C3").Value = 171.82536395782364
$ws.Range("D3").Value = 215.17433980092929
$ws.Range("E3").Value = 183.72897075034024

# Selection now only highlights the updated columns
$ws.Range("B1:E3").Select()
